$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "ログイン者(人)" (column C) and "入力者(人)" (column D) counts per
# university row, per the refreshed upload. Each tuple is
# (row, newLoginCount, newEntryCount); $null means that column is unchanged
# for that row.
$updates = @(
    @(2, 14, 11),
    @(3, 35, 24),
    @(4, 76, 50),
    @(5, 63, 53),
    @(6, 101, 72),
    @(7, 65, 39),
    @(8, 41, 26),
    @(9, 48, 37),
    @(10, 26, 19),
    @(11, 54, 40),
    @(12, 30, 23),
    @(13, 47, 33),
    @(14, 82, 64),
    @(15, 51, 39),
    @(16, 65, 47),
    @(17, 24, 17),
    @(18, 40, 27),
    @(19, 37, 26),
    @(20, 27, 20),
    @(21, 55, 37),
    @(22, $null, 20),
    @(23, 19, 11),
    @(24, 60, 44),
    @(25, 31, 23),
    @(26, 29, 18),
    @(27, 52, 36),
    @(28, 50, 34),
    @(29, 56, 36),
    @(30, 44, 29),
    @(31, 51, 36),
    @(32, 38, 21),
    @(33, 43, 28),
    @(34, 51, 31),
    @(35, 52, 39),
    @(36, 37, 20),
    @(37, 48, 28),
    @(38, 58, 41),
    @(39, 55, 47),
    @(40, 68, 54),
    @(41, 47, 33),
    @(42, 65, 40),
    @(43, 55, 41),
    @(44, 57, 46),
    @(45, 32, 18),
    @(46, 59, 40),
    @(47, 44, 34),
    @(48, 44, 25),
    @(49, 53, 42),
    @(50, 32, 23),
    @(51, 34, 25),
    @(52, 43, 27),
    @(53, 78, 62),
    @(54, 46, 22),
    @(55, 63, 41),
    @(56, 80, 59),
    @(57, 62, 33),
    @(58, 59, 46),
    @(59, 59, 45),
    @(60, 50, 30),
    @(61, 28, 19),
    @(62, 53, 34),
    @(63, 95, 80),
    @(64, 29, 20),
    @(65, 28, 23),
    @(66, 58, 36),
    @(67, 68, 49),
    @(68, 43, 28),
    @(69, 92, 68),
    @(70, 45, 34),
    @(71, 65, 53),
    @(72, 37, 15),
    @(73, 54, 32),
    @(74, 105, 85),
    @(75, 107, 104),
    @(76, 24, 11),
    @(77, 123, 118),
    @(78, 72, 52),
    @(79, 98, 97),
    @(80, 107, 86),
    @(81, 69, 43),
    @(82, 24, 13),
    @(83, 62, 48),
    @(84, 179, 98),
    @(86, $null, 7),
    @(92, 200, 122),
    @(93, 4798, 3412)
)

foreach ($u in $updates) {
    $row = $u[0]
    $newC = $u[1]
    $newD = $u[2]
    if ($null -ne $newC) {
        $ws.Cells.Item($row, 3).Value = $newC
    }
    if ($null -ne $newD) {
        $ws.Cells.Item($row, 4).Value = $newD
    }
}

$ws.Range("A1").Select()
